$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values regenerated to filter save games.
# Columns: B=TB, C=d2S, D=K, E=IP, F=Win, G=sum (F unchanged)
$data = @{
    2 = @(0.003994804209775715, 0.0000006633126561350622, 0.8054896365839992, 0.496779210170732, 1.306264314277163)
    3 = @(0.0008583669626518464, 0.002777888934908601, 0.8054896365839992, 8.660232485948974, 9.469358378430535)
    4 = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 8.660232485948974, 9.25842976318579)
    5 = @(0.003994804209775715, 225321.0684179339, 337.1190423067083, 616238.5361209477, 841896.7275759925)
    6 = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 645.3272768299601, 647.1208870175894)
    7 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
